$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("translations")

# New translation rows for the single-view button tooltips (print, json, watchlist add/remove)
$ws.Range("A50").Value = "print_entry"
$ws.Range("C50").Value = "imprimer"
$ws.Range("D50").Value = "print this entry"

$ws.Range("A51").Value = "json_data"
$ws.Range("C51").Value = "données JSON"
$ws.Range("D51").Value = "JSON data"

$ws.Range("A52").Value = "add_to_watchlist"
$ws.Range("C52").Value = "ajouter à la watchlist"
$ws.Range("D52").Value = "add to watchlist"

$ws.Range("A53").Value = "remove_from_watchlist"
$ws.Range("C53").Value = "supprimer de la watchlist"
$ws.Range("D53").Value = "remove from watchlist"

$ws.Range("C54").Select()
